$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for specific rows to repulled data
$ws.Range("F2").Value = -1
$ws.Range("F3").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("F12").Value = 3
$ws.Range("F16").Value = -1
$ws.Range("F18").Value = 3
$ws.Range("F19").Value = 4
$ws.Range("F20").Value = -1
$ws.Range("F27").Value = 1
$ws.Range("F30").Value = -2
$ws.Range("F33").Value = -2
$ws.Range("F36").Value = -1
$ws.Range("F37").Value = 1
